$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Instructions")
$ws3 = $wb.Worksheets.Item("ZEVs Supplied")

$ws1.Unprotect()
$ws1.Range("A1").Value = 'Please fill out the "ZEVs Supplied" sheet.'

$ws1.Range("A9").Value = 'With respect to the "ZEVs Supplied" sheet:'
$ws1.Range("A10").Value = '(1) Only the vehicles listed in the "Valid Vehicles" sheet may be used.'
$ws1.Range("A11").Value = "(2) VINs must be exactly 17 characters."
$ws1.Range("A13").Value = "(4) No more than 2000 records may be entered."
$ws1.Range("A12").Value = "(3) Dates must be of the YYYY-MM-DD format."

$ws1.Range("A3").Value = "With respect to this entire document:"
$ws1.Range("A4").Value = "(1) Please do not change the existing cell formatting."
$ws1.Range("A5").Value = "(2) Please do not add any cell formatting."
$ws1.Range("A6").Value = "(3) Please do not add any data validation."
$ws1.Range("A7").Value = "(4) Please do not change any of the sheet names or header names."

$ws1.Columns.Item(1).ColumnWidth = 59.5
$ws1.Protect()

$ws3.Range("E2:E2001").Validation.Delete()
$ws3.Range("E2:E1048576").NumberFormat = "@"
